$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix corporate/person name strings: comma -> period, drop other periods ---
$ws.Range("E30").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E35").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F35").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E57").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("F57").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("E66").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E67").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F67").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E70").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E99").Value = 'RICCOTTI. MARIANA EDITH'
$ws.Range("E100").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F100").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E129").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E130").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E131").Value = 'RICCOTTI. MARIANA EDITH'

# --- Fix "Importe" numeric-looking text: remove thousands separators (.), convert decimal comma to period ---
# Force column H to remain text (it was stored as text strings, not real numbers) before writing new values
$ws.Range("H2:H163").NumberFormat = "@"

$ws.Range("H2").Value = '1730.00'
$ws.Range("H3").Value = '766.00'
$ws.Range("H4").Value = '131000.10'
$ws.Range("H5").Value = '142.50'
$ws.Range("H6").Value = '3729.60'
$ws.Range("H7").Value = '40000.00'
$ws.Range("H8").Value = '2210.60'
$ws.Range("H9").Value = '4415.97'
$ws.Range("H10").Value = '23450.00'
$ws.Range("H11").Value = '131196.99'
$ws.Range("H12").Value = '1186.50'
$ws.Range("H13").Value = '18049.03'
$ws.Range("H14").Value = '5926.80'
$ws.Range("H15").Value = '2081.87'
$ws.Range("H16").Value = '5323.63'
$ws.Range("H17").Value = '762.40'
$ws.Range("H18").Value = '461.20'
$ws.Range("H19").Value = '6319.49'
$ws.Range("H20").Value = '690.00'
$ws.Range("H21").Value = '58.00'
$ws.Range("H22").Value = '1000.00'
$ws.Range("H23").Value = '1245.00'
$ws.Range("H24").Value = '8322.00'
$ws.Range("H25").Value = '795.00'
$ws.Range("H26").Value = '4010.00'
$ws.Range("H27").Value = '764.56'
$ws.Range("H28").Value = '8447.00'
$ws.Range("H29").Value = '1545.70'
$ws.Range("H30").Value = '225.00'
$ws.Range("H31").Value = '12956.50'
$ws.Range("H32").Value = '13079.35'
$ws.Range("H33").Value = '2887.00'
$ws.Range("H34").Value = '0.14'
$ws.Range("H35").Value = '165.56'
$ws.Range("H36").Value = '10461.81'
$ws.Range("H37").Value = '11.00'
$ws.Range("H38").Value = '629.28'
$ws.Range("H39").Value = '20.90'
$ws.Range("H40").Value = '53.98'
$ws.Range("H41").Value = '26554.88'
$ws.Range("H42").Value = '1280.00'
$ws.Range("H43").Value = '50164.32'
$ws.Range("H44").Value = '5642.00'
$ws.Range("H45").Value = '28.34'
$ws.Range("H46").Value = '1242.19'
$ws.Range("H47").Value = '4575.14'
$ws.Range("H48").Value = '653.64'
$ws.Range("H49").Value = '640.00'
$ws.Range("H50").Value = '950.00'
$ws.Range("H51").Value = '2937.00'
$ws.Range("H52").Value = '3832.57'
$ws.Range("H53").Value = '888.00'
$ws.Range("H54").Value = '1746.29'
$ws.Range("H55").Value = '250.00'
$ws.Range("H56").Value = '72.00'
$ws.Range("H57").Value = '55.00'
$ws.Range("H58").Value = '120.00'
$ws.Range("H59").Value = '13674.00'
$ws.Range("H60").Value = '155.00'
$ws.Range("H61").Value = '11355.90'
$ws.Range("H62").Value = '400.00'
$ws.Range("H63").Value = '350.00'
$ws.Range("H64").Value = '32941.00'
$ws.Range("H65").Value = '5211.80'
$ws.Range("H66").Value = '4515.60'
$ws.Range("H67").Value = '2142.71'
$ws.Range("H68").Value = '5345.48'
$ws.Range("H69").Value = '180.00'
$ws.Range("H70").Value = '7092.00'
$ws.Range("H71").Value = '304.50'
$ws.Range("H72").Value = '1337.07'
$ws.Range("H73").Value = '955.20'
$ws.Range("H74").Value = '1740.00'
$ws.Range("H75").Value = '36635.00'
$ws.Range("H76").Value = '940.00'
$ws.Range("H77").Value = '91770.00'
$ws.Range("H78").Value = '12998.31'
$ws.Range("H79").Value = '0.24'
$ws.Range("H80").Value = '75.00'
$ws.Range("H81").Value = '4349.47'
$ws.Range("H82").Value = '2228.80'
$ws.Range("H83").Value = '1753.00'
$ws.Range("H84").Value = '809.50'
$ws.Range("H85").Value = '956.26'
$ws.Range("H86").Value = '7780.00'
$ws.Range("H87").Value = '86.50'
$ws.Range("H88").Value = '836.50'
$ws.Range("H89").Value = '3131.10'
$ws.Range("H90").Value = '7.48'
$ws.Range("H91").Value = '670.00'
$ws.Range("H92").Value = '287.00'
$ws.Range("H93").Value = '8839.50'
$ws.Range("H94").Value = '73.70'
$ws.Range("H95").Value = '378.00'
$ws.Range("H96").Value = '142.34'
$ws.Range("H97").Value = '336.00'
$ws.Range("H98").Value = '27.83'
$ws.Range("H99").Value = '1600.00'
$ws.Range("H100").Value = '832.00'
$ws.Range("H101").Value = '498.00'
$ws.Range("H102").Value = '875.00'
$ws.Range("H103").Value = '145.20'
$ws.Range("H104").Value = '390.00'
$ws.Range("H105").Value = '2755.73'
$ws.Range("H106").Value = '794487.28'
$ws.Range("H107").Value = '870.00'
$ws.Range("H108").Value = '5831.00'
$ws.Range("H109").Value = '200.00'
$ws.Range("H110").Value = '500.00'
$ws.Range("H111").Value = '1400.00'
$ws.Range("H112").Value = '150.00'
$ws.Range("H113").Value = '11043.92'
$ws.Range("H114").Value = '290.00'
$ws.Range("H115").Value = '9317.00'
$ws.Range("H116").Value = '265.20'
$ws.Range("H117").Value = '17545.00'
$ws.Range("H118").Value = '250.00'
$ws.Range("H119").Value = '2100.00'
$ws.Range("H120").Value = '1210.00'
$ws.Range("H121").Value = '2238.00'
$ws.Range("H122").Value = '689.52'
$ws.Range("H123").Value = '1800.00'
$ws.Range("H124").Value = '450.00'
$ws.Range("H125").Value = '750.00'
$ws.Range("H126").Value = '120.00'
$ws.Range("H127").Value = '9539.01'
$ws.Range("H128").Value = '198.00'
$ws.Range("H129").Value = '38.00'
$ws.Range("H130").Value = '60.00'
$ws.Range("H131").Value = '3800.00'
$ws.Range("H132").Value = '363.00'
$ws.Range("H133").Value = '3325.00'
$ws.Range("H134").Value = '130.00'
$ws.Range("H135").Value = '24.40'
$ws.Range("H136").Value = '606.96'
$ws.Range("H137").Value = '2235.81'
$ws.Range("H138").Value = '1745.00'
$ws.Range("H139").Value = '99.00'
$ws.Range("H140").Value = '9279.10'
$ws.Range("H141").Value = '675.00'
$ws.Range("H142").Value = '888.80'
$ws.Range("H143").Value = '3919.00'
$ws.Range("H144").Value = '231.84'
$ws.Range("H145").Value = '29.00'
$ws.Range("H146").Value = '2410.50'
$ws.Range("H147").Value = '25319.50'
$ws.Range("H148").Value = '2076.26'
$ws.Range("H149").Value = '1831.00'
$ws.Range("H150").Value = '2532.53'
$ws.Range("H151").Value = '64405.71'
$ws.Range("H152").Value = '2738.29'
$ws.Range("H153").Value = '2200.00'
$ws.Range("H154").Value = '31150.00'
$ws.Range("H155").Value = '1194.92'
$ws.Range("H156").Value = '823.00'
$ws.Range("H157").Value = '316948.14'
$ws.Range("H158").Value = '33400.00'
$ws.Range("H159").Value = '227379.37'
$ws.Range("H160").Value = '41000.00'
$ws.Range("H161").Value = '1129.26'
$ws.Range("H162").Value = '4000.00'
$ws.Range("H163").Value = '3800.00'
